$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new data row for the 2022-Q3 quarter
#    right after the header, pushing the existing quarter rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
# Inserting a row copies the formatting of the row above (the bold header),
# so clear direct formatting on the new row's text/number cells first ...
$summary.Range("B2:D2").Style = "Normal"
# ... then restore the plain numbered-index style used by the other rows'
# first column (matches what A3:A7 already use).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 2.16

# ---------------------------------------------------------------------------
# 2) Add a brand-new "2022-Q3" sheet, positioned right after "总计" and
#    before "2022-Q2" (all the later quarter sheets shift back by one slot).
#    Cloning the "2022-Q2" sheet keeps identical layout/formatting.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

function Set-TextCell($range, [string]$text) {
    # A leading apostrophe forces the value to be stored as text even when
    # it looks numeric; resetting the style afterwards drops the implicit
    # "quote prefix" direct formatting Excel applies so the cell keeps the
    # same plain style as its neighbours.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextCell $q3.Range("D2") "20.44"
Set-TextCell $q3.Range("E2") "90.19"
Set-TextCell $q3.Range("F2") "6.09"
Set-TextCell $q3.Range("G2") "1.2448"
$q3.Range("H2").Value = 6

Set-TextCell $q3.Range("D3") "15.02"
Set-TextCell $q3.Range("E3") "90.19"
Set-TextCell $q3.Range("F3") "6.09"
Set-TextCell $q3.Range("G3") "0.9147"
$q3.Range("H3").Value = 6

# Restore the originally active sheet/selection (the summary sheet was the
# active tab before this edit).
$summary.Select()
$summary.Range("A1").Select() | Out-Null
